$wb = $excel.ActiveWorkbook

# --- Sheet: All Orders ---
$wsOrders = $wb.Worksheets.Item("All Orders")
$wsOrders.Range("H3").Value = "DELIVERED"
$wsOrders.Range("I3").Value = "PAID"

# --- Sheet: Daily Summary ---
$wsSummary = $wb.Worksheets.Item("Daily Summary")
$wsSummary.Range("C2").Value = 2
$wsSummary.Range("F2").Value = 240
$wsSummary.Range("G2").Value = 285
